$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
$lastCol = $usedRange.Column + $usedRange.Columns.Count - 1

# Locate the "Recorded By" column dynamically from the header row
$recordedByCol = 7
for ($c = 1; $c -le $lastCol; $c++) {
    $header = $ws.Cells.Item(1, $c).Value2
    if ($header -eq "Recorded By") {
        $recordedByCol = $c
    }
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $val = $cell.Value2
    if ($null -ne $val -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ","
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }
        $reversed = $trimmed[($trimmed.Length - 1)..0]
        $newVal = [string]::Join(", ", $reversed)
        $cell.Value = $newVal
    }
}
